$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Fill in the new "Purposes" column (C) for the existing rows,
#    fix up a couple of other cell values, and add the two new rows
#    (jumper wires) to the bill of materials.
# -----------------------------------------------------------------

# New rows for jumper wires (written B14 first, then B13, matching
# the order the shared-string table was populated by the author)
$ws.Range("B14").Value = "Male-to-male jumper wires"
$ws.Range("B13").Value = "Female-to-male jumper wires"

$ws.Range("C5").Value = "Microcontroller board, controls the program"
$ws.Range("C6").Value = "Used to display menus and options to a user"
$ws.Range("C7").Value = "Control current flow to LEDs"
$ws.Range("C9").Value = "Allows user to interact with menu options"
$ws.Range("C10").Value = "Prototyping board, used to design board layout"
$ws.Range("C11").Value = "Detects when board is uneven"
$ws.Range("C12").Value = "Provides audio feedback to user"
$ws.Range("C8").Value = "Provide visual feedback to user"

$ws.Range("C13").Value = "Connects components"
$ws.Range("C14").Value = "Connects components"

# Tilt sensor quantity bumped from 2 to 4
$ws.Range("A11").Value = 4
